$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "30.582.51"
$ws.Range("E2").Value = "  +0.25%  "

Set-TextValue "D3" "2.114.79"
$ws.Range("E3").Value = "  +0.17%  "

Set-TextValue "D5" "345.95"
$ws.Range("E5").Value = "  +3.35%  "

Set-TextValue "D6" "1.011"
$ws.Range("E6").Value = "  +0.92%  "

Set-TextValue "D7" "0.5253"
$ws.Range("E7").Value = "  -0.27%  "

Set-TextValue "D8" "0.4513"
$ws.Range("E8").Value = "  -1.14%  "

Set-TextValue "D9" "53.75"
$ws.Range("E9").Value = "  -0.31%  "

Set-TextValue "D10" "0.09013"
$ws.Range("E10").Value = "  +0.10%  "

Set-TextValue "D11" "1.171"
$ws.Range("E11").Value = "  -1.32%  "

Set-TextValue "D12" "24.44"
$ws.Range("E12").Value = "  -0.33%  "

Set-TextValue "D13" "2.117.18"
$ws.Range("E13").Value = "  +1.15%  "

Set-TextValue "D14" "6.806"
$ws.Range("E14").Value = "  -0.32%  "

Set-TextValue "D15" "8.082"
$ws.Range("E15").Value = "  +2.98%  "

Set-TextValue "D16" "99.41"
$ws.Range("E16").Value = "  +2.51%  "

Set-TextValue "D17" "0.00001171"
$ws.Range("E17").Value = "  +3.28%  "

Set-TextValue "D18" "1.013"
$ws.Range("E18").Value = "  +0.95%  "

Set-TextValue "D19" "0.06709"
$ws.Range("E19").Value = "  +1.30%  "

Set-TextValue "D20" "19.36"
$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("E21").Value = "  +1.03%  "

Set-TextValue "D22" "6.331"
$ws.Range("E22").Value = "  +0.17%  "

Set-TextValue "D23" "30.674.65"
$ws.Range("E23").Value = "  +0.36%  "

Set-TextValue "D24" "12.74"
$ws.Range("E24").Value = "  +2.92%  "

Set-TextValue "D25" "2.384"
$ws.Range("E25").Value = "  +1.15%  "

Set-TextValue "D26" "2.363.70"
$ws.Range("E26").Value = "  +1.02%  "

Set-TextValue "D27" "22.40"
$ws.Range("E27").Value = "  -0.15%  "

Set-TextValue "D28" "165.51"
$ws.Range("E28").Value = "  +1.00%  "

Set-TextValue "D29" "2.537"
$ws.Range("E29").Value = "  -1.81%  "

Set-TextValue "D30" "135.02"
$ws.Range("E30").Value = "  +1.45%  "

Set-TextValue "D31" "1.194"
$ws.Range("E31").Value = "  -0.59%  "

Set-TextValue "D32" "0.1073"
$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("E33").Value = "  -2.41%  "

Set-TextValue "D34" "6.364"
$ws.Range("E34").Value = "  +3.14%  "

Set-TextValue "D35" "3.963"
$ws.Range("E35").Value = "  +0.83%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D36" "5.899"
$ws.Range("E36").Value = "  +5.44%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D37" "10.20"
$ws.Range("E37").Value = "  -2.56%  "

Set-TextValue "D38" "0.02635"
$ws.Range("E38").Value = "  +1.93%  "

Set-TextValue "D39" "0.06841"
$ws.Range("E39").Value = "  -0.06%  "

Set-TextValue "D40" "0.2330"
$ws.Range("E40").Value = "  +1.36%  "

Set-TextValue "D41" "12.65"
$ws.Range("E41").Value = "  -1.06%  "

Set-TextValue "D42" "0.6871"
$ws.Range("E42").Value = "  -0.86%  "

Set-TextValue "D43" "1.265"
$ws.Range("E43").Value = "  +1.26%  "

Set-TextValue "D44" "14.83"
$ws.Range("E44").Value = "  +5.30%  "

Set-TextValue "D45" "0.6424"
$ws.Range("E45").Value = "  +0.28%  "

Set-TextValue "D46" "2.314"
$ws.Range("E46").Value = "  -2.31%  "

Set-TextValue "D47" "3.738"
$ws.Range("E47").Value = "  +2.18%  "

Set-TextValue "D48" "0.00000000360"
$ws.Range("E48").Value = "  +1.37%  "

Set-TextValue "D49" "1.254"
$ws.Range("E49").Value = "  +0.15%  "

Set-TextValue "D50" "82.86"
$ws.Range("E50").Value = "  -0.94%  "

Set-TextValue "D51" "0.07287"
$ws.Range("E51").Value = "  +2.68%  "
